$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (B1) and short name (B2) on the input sheet
$wsInput.Range("B1").Value = "4300-MS-EI-DB-SAR-REC-RNI-INT-FFConMONTHLYonLASTSUNDAY-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1st"
$wsInput.Range("B2").Value = "430q"

# Mirror the product name on the output sheet
$wsOutput.Range("B1").Value = "4300-MS-EI-DB-SAR-REC-RNI-INT-FFConMONTHLYonLASTSUNDAY-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1st"

# Move the selection on the input sheet
$wsInput.Range("B3").Select()

# Make the output sheet the active tab
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
